$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain plain text so values like "1.001" or "309.96"
# are not auto-converted to floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.853.85'
$ws.Range("E2").Value = '  -1.71%  '
$ws.Range("D3").Value = '1.804.66'
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '309.96'
$ws.Range("E5").Value = '  -1.40%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").Value = '0.4656'
$ws.Range("E7").Value = '  +3.82%  '
$ws.Range("D8").Value = '0.3712'
$ws.Range("E8").Value = '  -2.04%  '
$ws.Range("D9").Value = '0.07383'
$ws.Range("E9").Value = '  -0.76%  '
$ws.Range("D10").Value = '0.8723'
$ws.Range("E10").Value = '  -1.59%  '
$ws.Range("E11").Value = '  -2.78%  '
$ws.Range("D12").Value = '1.760.33'
$ws.Range("E12").Value = '  -3.54%  '
$ws.Range("E13").Value = '  -1.54%  '
$ws.Range("D14").Value = '92.65'
$ws.Range("E14").Value = '  -1.02%  '
$ws.Range("D15").Value = '6.484'
$ws.Range("E15").Value = '  -3.74%  '
$ws.Range("D16").Value = '0.07027'
$ws.Range("E16").Value = '  -1.43%  '
$ws.Range("D17").Value = '1.001'
$ws.Range("D18").Value = '0.000008714'
$ws.Range("E18").Value = '  -1.09%  '
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("E20").Value = '  -2.99%  '
$ws.Range("D21").Value = '26.855.72'
$ws.Range("E21").Value = '  -1.72%  '
$ws.Range("D22").Value = '5.304'
$ws.Range("E22").Value = '  -1.61%  '
$ws.Range("E23").Value = '  -3.11%  '
$ws.Range("D24").Value = '2.015.15'
$ws.Range("E24").Value = '  -1.79%  '
$ws.Range("D25").Value = '1.892'
$ws.Range("E25").Value = '  -3.87%  '
$ws.Range("D26").Value = '151.47'
$ws.Range("E26").Value = '  -0.15%  '
$ws.Range("D27").Value = '18.32'
$ws.Range("E27").Value = '  -1.74%  '
$ws.Range("D28").Value = '2.150'
$ws.Range("E28").Value = '  -7.05%  '
$ws.Range("D29").Value = '5.273'
$ws.Range("E29").Value = '  -2.20%  '
$ws.Range("D30").Value = '115.80'
$ws.Range("E30").Value = '  -1.70%  '
$ws.Range("D31").Value = '0.08945'
$ws.Range("E31").Value = '  +0.53%  '
$ws.Range("D32").Value = '0.7586'
$ws.Range("E32").Value = '  -4.38%  '
$ws.Range("E33").Value = '  -3.97%  '
$ws.Range("D34").Value = '4.455'
$ws.Range("E34").Value = '  -3.21%  '
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("D36").Value = '1.000'
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("D37").Value = '1.105'
$ws.Range("E37").Value = '  -0.60%  '
$ws.Range("D38").Value = '0.01960'
$ws.Range("E38").Value = '  -1.11%  '
$ws.Range("D39").Value = '0.05257'
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '2.930'
$ws.Range("E40").Value = '  +2.09%  '
$ws.Range("D41").Value = '7.247'
$ws.Range("E41").Value = '  -1.41%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = '2.382'
$ws.Range("E42").Value = '  +1.73%  '
$ws.Range("E43").Value = '  -0.89%  '
$ws.Range("D44").Value = '0.1663'
$ws.Range("E44").Value = '  -3.07%  '
$ws.Range("D45").Value = '8.509'
$ws.Range("E45").Value = '  -1.93%  '
$ws.Range("D46").Value = '0.5002'
$ws.Range("E46").Value = '  -1.32%  '
$ws.Range("D47").Value = '10.32'
$ws.Range("E47").Value = '  -3.12%  '
$ws.Range("D48").Value = '104.01'
$ws.Range("E48").Value = '  -1.28%  '
$ws.Range("D49").Value = '1.000'
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("E50").Value = '  -1.91%  '
$ws.Range("D51").Value = '0.06294'
$ws.Range("E51").Value = '  -1.77%  '
